$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 657
$ws.Range("F5").Value = 2954
$ws.Range("F7").Value = 241
$ws.Range("F10").Value = 6936
$ws.Range("F11").Value = 42
$ws.Range("F12").Value = 89
$ws.Range("F13").Value = 355
$ws.Range("F15").Value = 1498
$ws.Range("F16").Value = 1121
$ws.Range("F17").Value = 2247
$ws.Range("F18").Value = 1498
$ws.Range("F20").Value = 123
$ws.Range("F21").Value = 1117
$ws.Range("F22").Value = 132
$ws.Range("F23").Value = 181
$ws.Range("F26").Value = 1740
$ws.Range("F28").Value = 1034
$ws.Range("F29").Value = 37
$ws.Range("F30").Value = 1668
$ws.Range("F31").Value = 1227
$ws.Range("F33").Value = 589
$ws.Range("F34").Value = 38
$ws.Range("F36").Value = 18
$ws.Range("F37").Value = 2489
$ws.Range("F40").Value = 12
$ws.Range("F42").Value = 18
$ws.Range("G42").Value = 68
$ws.Range("F44").Value = 320
$ws.Range("F47").Value = 159
$ws.Range("F49").Value = 415

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 213
$ws.Range("F15").Value = 58
$ws.Range("F19").Value = 46
$ws.Range("F20").Value = 55
$ws.Range("F23").Value = 478
$ws.Range("F37").Value = 40

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1690
$ws.Range("F7").Value = 1855
$ws.Range("F8").Value = 2738
$ws.Range("F9").Value = 1023
$ws.Range("F10").Value = 938
$ws.Range("F12").Value = 274
$ws.Range("F13").Value = 1484
$ws.Range("F14").Value = 7372

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 657
$ws.Range("F5").Value = 2954
$ws.Range("F6").Value = 1690
$ws.Range("F8").Value = 2738
$ws.Range("F9").Value = 6936
$ws.Range("F10").Value = 1023
$ws.Range("F11").Value = 42
$ws.Range("F12").Value = 355
$ws.Range("F14").Value = 274
$ws.Range("F15").Value = 1121
$ws.Range("F16").Value = 2247
$ws.Range("F17").Value = 1498
$ws.Range("F18").Value = 123
$ws.Range("F20").Value = 1117
$ws.Range("F22").Value = 1740
$ws.Range("F24").Value = 37
$ws.Range("F25").Value = 1668
$ws.Range("F26").Value = 1227
$ws.Range("F29").Value = 589
$ws.Range("F30").Value = 38
$ws.Range("F31").Value = 55
$ws.Range("F34").Value = 478
$ws.Range("F37").Value = 18
$ws.Range("F38").Value = 2489
$ws.Range("F42").Value = 18
$ws.Range("G42").Value = 68
$ws.Range("F44").Value = 320
$ws.Range("F48").Value = 415

Write-Output "Updated 72 cells across 4 sheets"
